$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# --- Add hyperlinks to existing URL text in J5 and J6 (keep original formatting) ---
$ws.Range("J5").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Hyperlinks.Add($ws.Range("J5"), $ws.Range("J5").Value2) | Out-Null
$ws.Range("Z1").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("Z1").Clear()

$ws.Range("J6").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Hyperlinks.Add($ws.Range("J6"), $ws.Range("J6").Value2) | Out-Null
$ws.Range("Z1").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("Z1").Clear()

# --- Fill in row 13: Load cell + amplifier (DEBO HX711-01) ---
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Joy-It"
$ws.Range("D13").WrapText = $true
$ws.Range("D13").Value = "Load cell + amplifier "
$ws.Range("G13").Value = "DEBO HX711-01"
$ws.Range("E13").Value = "SEN-HX711-01"
$ws.Range("F13").Value = "Reichelt"
$ws.Range("H13").Value = 6.4
$ws.Range("J13").Value = "https://www.reichelt.com/de/en/shop/product/developer_boards_-_a_d_converter_-_balance-316296#closemodal"

$ws.Rows.Item(13).RowHeight = 34.5

$ws.Range("H17").Select()

Write-Output "done"
